# Update countries & provincias Spain
# Applies the "9 Abril 2020 08:52" data refresh to the Pais sheet:
#  - refreshed timestamp banner
#  - Ucrania, Bosnia y Herzegovina and Paraguay each picked up enough new
#    cases overnight to overtake their neighbours in the (descending,
#    sorted-by-Casos-totales) table, pushing those rows down by one
#  - a few other rows (Chequia, Rumania, El Salvador) just got refreshed
#    figures without changing rank

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $country, $values) {
    if ($country -ne $null) {
        $ws.Cells.Item($row, 1).Value = $country
    }
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($row, 2 + $i).Value = $values[$i]
    }
}

# --- Header: refresh timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 9 de Abril de 2020 a las 08:52"

# --- Minor data-only updates (no re-ranking) ---
Set-Row 30 $null @(5335, 23, 243, 4988, 96, 5, 104)   # Chequia
Set-Row 32 $null @(4761, 0, 528, 4006, 162, 7, 227)   # Rumania
$ws.Cells.Item(130, 6).Value = 3                      # El Salvador, Nuevos casos

# --- Ucrania overtakes Argentina/Sudafrica/Grecia (rows 51-54) ---
Set-Row 51 "Ucrania"   @(1892, 224, 45, 1790, 33, 5, 57)
Set-Row 52 "Grecia"    @(1884, 0, 269, 1532, 84, 0, 83)
Set-Row 53 "Sudafrica" @(1845, 0, 95, 1732, 7, 0, 18)
Set-Row 54 "Argentina" @(1795, 0, 358, 1372, 96, 0, 65)

# --- Bosnia y Herzegovina overtakes Barein/Azerbaiyan (rows 72-74) ---
Set-Row 72 "Bosnia y Herzegovina" @(823, 19, 95, 693, 4, 1, 35)
Set-Row 73 "Barein"               @(823, 0, 477, 341, 3, 0, 5)
Set-Row 74 "Azerbaiyan"           @(822, 0, 63, 751, 23, 0, 8)

# --- Paraguay overtakes Gibraltar (rows 125-126) ---
Set-Row 125 "Paraguay"  @(124, 5, 18, 101, 1, 0, 5)
Set-Row 126 "Gibraltar" @(120, 0, 60, 60, 0, 0, 0)
